# Update from MV -datos-: append three new daily rows (06-11-2021,
# 07-11-2021, 08-11-2021) to Sheet1, mirroring the existing rows
# (same B/C values as the prior last row, 322 / 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDates = @("06-11-2021", "07-11-2021", "08-11-2021")
$startRow = 284

for ($i = 0; $i -lt $newDates.Count; $i++) {
    $row = $startRow + $i

    # Force column A to be read as literal text so Excel doesn't
    # reinterpret the dd-mm-yyyy string as a date serial number.
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $newDates[$i]
    # Drop the temporary text format again so the cell ends up with
    # the same (default/no explicit style) formatting as the rest of
    # the data rows in column A.
    $ws.Range("A$row").ClearFormats()

    $ws.Range("B$row").Value = 322
    $ws.Range("C$row").Value = 0
}
